$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill empty cells in the SIZE (cm) column (K) with "-" for rows where the
# measurement was not applicable / not recorded, matching column A formatting.
$rows = @(5,7,12,13,14,22,41,42,43,49,50,53,56,59,60,63,64,65,71,75,76,77,78,80,88,92,96,105,106,107,108,115,116,118,119,122,123,126,128,129,131,136,138)
foreach ($r in $rows) {
    $src = $ws.Range("A" + $r)
    $dst = $ws.Range("K" + $r)
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $dst.Value = "-"
}

# Restore cursor / selection to where the author left it after editing.
$ws.Range("I133").Select()
